# Add season record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - match style of existing header cells (A1:AC1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerSrc = $ws.Range("AC1")
$headerDst = $ws.Range("AD1:AF1")
$headerDst.Font.Bold = $headerSrc.Font.Bold
$headerDst.HorizontalAlignment = $headerSrc.HorizontalAlignment
$headerDst.VerticalAlignment = $headerSrc.VerticalAlignment
$headerDst.Borders.LineStyle = $headerSrc.Borders.LineStyle

# Data rows 2-48: every player on this roster shares the team's season record.
for ($row = 2; $row -le 48; $row++) {
    $ws.Cells.Item($row, 30).Value = 68
    $ws.Cells.Item($row, 31).Value = 94
    $ws.Cells.Item($row, 32).Value = 0
}
